$wb = $excel.ActiveWorkbook

# The file "9a943c75-0ca4-4009-84b3-9fdfbe3095e3.md" has moved from
# "Ready for handoff" to "In Translation" status. Update the Status
# cell for that file's row on the Overview sheet as well as the
# per-language (zh-cn / de-de) detail sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "In Translation"
$overview.Range("C2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "In Translation"
